# endpointassay.xlsx: add "UsedEnzyme" column (new X), shift the old
# "Comment" column content into a brand-new Y column, and append a
# "(nm)" unit hint to the Wavelength format-help cell (J4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column J (Wavelength) : clarify the format-help text with its unit ---
$ws.Range("J4").Value = "# format: nombre entier, ne pas spécifier d'unité (nm)"

# --- Column X becomes the new "UsedEnzyme" column -------------------------
$ws.Range("X1").Value = "UsedEnzyme"
$ws.Range("X2").Value = "# Enzyme utilisée"
$ws.Range("X3").Value = "#string"
$ws.Range("X4").Value = "# format: texte"
$ws.Range("X5").Value = "# ex:"

# --- Column Y is new: it carries what used to be the "Comment" column ----
$ws.Range("Y1").Value = "Comment"
$ws.Range("Y2").Value = "# Commentaire"
$ws.Range("Y3").Value = "#string"
$ws.Range("Y4").Value = "# format: texte libre"
# Y5 stays blank (the original Comment-example cell was empty too).
